# Replace the open/close/high/low price, shares_outstanding, and fixed_ticker
# columns (D:I) for every data row (2-43) with VRSK-specific values, undoing
# the accidental mix-in of other companies price/ticker data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row, open_price, close_price, high_price, low_price
$data = @(
    ,@(2, 68.39171792802742, 71.8682861328125, 74.78936743670167, 67.98946450952722)
    ,@(3, 70.02944661163185, 74.80853271484375, 76.50371707975654, 68.84185231157937)
    ,@(4, 74.71276876524381, 68.5832748413086, 78.53412217642621, 66.74442885620535)
    ,@(5, 72.17477896840994, 69.9145278930664, 72.38548145445593, 66.17936672232582)
    ,@(6, 76.27386163751362, 74.30093383789062, 78.14143821415446, 73.51559300733864)
    ,@(7, 77.65299279842692, 81.67546844482422, 82.36503736864613, 77.17412630209759)
    ,@(8, 77.50932188335969, 78.10311889648438, 79.0991546566158, 75.2778025304494)
    ,@(9, 78.14141891593357, 79.14704132080078, 79.51097716954463, 76.81017085275688)
    ,@(10, 77.63386497167446, 79.30989837646484, 80.15270842216336, 75.36404155133368)
    ,@(11, 81.2828259120891, 83.57180786132812, 84.28052844091499, 79.43440039377839)
    ,@(12, 79.87495413761141, 81.45521545410156, 81.80957207851147, 79.50143875641325)
    ,@(13, 92.4020833516721, 95.8211898803711, 96.29047700593917, 91.23364791881356)
    ,@(14, 98.9433948983168, 101.9506759643555, 103.9906519515384, 96.07977524512415)
    ,@(15, 102.4678593897214, 105.9444351196289, 110.2350751929238, 102.0464618064929)
    ,@(16, 116.3837419163118, 114.7747497558594, 116.805146886669, 106.8064093676483)
    ,@(17, 102.5253191452748, 112.4474334716797, 113.4722074751345, 100.4853432072474)
    ,@(18, 128.2692524920409, 135.4374389648438, 135.5909779071337, 127.3576388158984)
    ,@(19, 141.8024910440722, 145.8397064208984, 147.8583214430095, 140.8604628778093)
    ,@(20, 152.4392708363372, 139.3163757324219, 153.6331414295598, 138.9505237742231)
    ,@(21, 143.9711234500206, 156.6920471191406, 160.0482830421848, 143.6046327590882)
    ,@(22, 129.6775533413751, 147.6577301025391, 154.5754156067351, 128.9819193385855)
    ,@(23, 163.7744820752355, 182.6265106201172, 182.6652146317478, 163.0486600303411)
    ,@(24, 180.0101331161973, 172.4891204833984, 189.0915785147856, 169.4361358623619)
    ,@(25, 201.4870052026782, 178.0966491699219, 203.49604417086, 177.1455104490742)
    ,@(26, 172.047045920862, 182.9646453857422, 184.5104111713112, 171.8428811782775)
    ,@(27, 170.9446979728797, 184.9677429199219, 185.2501447463326, 170.623349806715)
    ,@(28, 196.58108583235, 205.0555877685547, 210.6337405650339, 192.1634249991105)
    ,@(29, 222.2959310422345, 191.5088348388672, 222.5205076633652, 181.2171675499394)
    ,@(30, 210.8498366228578, 199.5820770263672, 217.2466287045212, 199.0832345524401)
    ,@(31, 169.292055488416, 186.4417572021484, 187.1375516797238, 167.7142823326185)
    ,@(32, 168.6417624693048, 179.468994140625, 183.4249084876406, 160.6317661078988)
    ,@(33, 174.4315662213555, 178.7481079101562, 184.2839061250982, 173.0746679137486)
    ,@(34, 187.4088207450072, 191.2212219238281, 194.6001688538694, 184.1776354762748)
    ,@(35, 221.2492693216344, 225.8865814208984, 235.6348145981541, 217.7367371448296)
    ,@(36, 232.4935291129258, 224.638687133789, 246.0196734084931, 219.2341573592451)
    ,@(37, 235.2971086767513, 238.9677429199219, 245.0821906275168, 227.0950441680811)
    ,@(38, 233.317600803873, 216.0043640136719, 233.317600803873, 215.3899166709513)
    ,@(39, 267.5293270643612, 259.7878112792969, 284.9775569097218, 258.1700256025808)
    ,@(40, 267.6587773312833, 273.0558776855469, 278.3038927352457, 259.548203103073)
    ,@(41, 274.6837698904532, 286.0899353027344, 288.498551877741, 264.9397437025851)
    ,@(42, 295.8857963941511, 295.4970092773437, 305.7147484398095, 272.8983778438283)
    ,@(43, 310.8056629826289, 278.2316589355469, 314.2597162708633, 269.8959834250493)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = 139714971 # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "VRSK"    # I: fixed_ticker
}

Write-Host "Updated $($data.Count) rows"